$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.072.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.959.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4886"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2955"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06919"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "107.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.961.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07781"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.470"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7004"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.088.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007733"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.226.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.518"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.507"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.846"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.201"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1052"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.392"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.441"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04934"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7558"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.735"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02013"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.707"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.528"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.122"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9067"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4458"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.192"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.029.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.406"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1258"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
